# "Add trustLevel in flow"
# The "Trust Boundaries" table lists two Enabler rows. Row 2 was labelled
# "Enablers2" and row 3 "Enablers" - swap them so row 2 reads "Enablers"
# and row 3 reads "Enablers2" (the other columns for both rows already
# share the same description/limit/level values).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trust Boundaries")

$ws.Range("A2").Value = "Enablers"
$ws.Range("A3").Value = "Enablers2"
